$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 121, shifting rows 121-179 down to 122-180.
$ws.Rows.Item(121).Insert(1)

# Populate the newly inserted row 121 with the new record.
$ws.Cells.Item(121, 1).Value = 4
$ws.Cells.Item(121, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(121, 3).Value = "Los Lagos"
$ws.Cells.Item(121, 4).Value = 44553
$ws.Cells.Item(121, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(121, 5).Value = 10
$ws.Cells.Item(121, 6).Value = "Fruta"
$ws.Cells.Item(121, 7).Value = 100108
$ws.Cells.Item(121, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(121, 9).Value = 100108005
$ws.Cells.Item(121, 10).Value = "Piña"
$ws.Cells.Item(121, 11).Value = "Caramelo"
$ws.Cells.Item(121, 12).Value = "Tercera"
$ws.Cells.Item(121, 13).Value = 200
$ws.Cells.Item(121, 14).Value = 19000
$ws.Cells.Item(121, 15).Value = 20000
$ws.Cells.Item(121, 16).Value = 19500
$ws.Cells.Item(121, 17).Value = "`$/caja 16 unidades"
$ws.Cells.Item(121, 18).Value = "Ecuador"
$ws.Cells.Item(121, 19).Value = 1219
$ws.Cells.Item(121, 20).Value = 16

$wb.Save()
